# Logs workbook: add a header row (Id / Description / Message) above the
# existing log rows, right-align the header, widen the Message column, and
# leave the selection on D12 - matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push all existing data down one row and create the new header row.
$ws.Rows("1:1").Insert()

$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "Message"
$ws.Range("A1:C1").HorizontalAlignment = -4152   # xlRight

# Column C (Message) needs to be noticeably wider now that it has a header.
$ws.Columns("C").ColumnWidth = 58.60807291666667

# Match the author's last on-screen selection.
$ws.Range("D12").Select()
